# Regenerate the "K" column (G) values for the save_data sheet.
# The underlying K (strikeout) values were recomputed (std/mean regen +
# calc/write of s_vals) and the resulting integers are written back into
# column G for every data row (rows 2-78).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2;  3  = 1;  4  = 1;  5  = 1;  6  = 0;  7  = 0;  8  = 0;  9  = 2;
    10 = 1;  11 = 2;  12 = 2;  13 = 3;  14 = 0;  15 = 0;  16 = 1;  17 = 1;
    18 = 2;  19 = 0;  20 = 3;  21 = 3;  22 = 1;  23 = 1;  24 = 1;  25 = 2;
    26 = 0;  27 = 1;  28 = 1;  29 = 2;  30 = 1;  31 = 0;  32 = 1;  33 = 1;
    34 = 0;  35 = 0;  36 = 0;  37 = 1;  38 = 2;  39 = 2;  40 = 1;  41 = 0;
    42 = 1;  43 = 2;  44 = 0;  45 = 0;  46 = 2;  47 = 1;  48 = 1;  49 = 1;
    50 = 0;  51 = 2;  52 = 2;  53 = 2;  54 = 1;  55 = 1;  56 = 1;  57 = 0;
    58 = 2;  59 = 0;  60 = 1;  61 = 1;  62 = 2;  63 = 2;  64 = 0;  65 = 1;
    66 = 1;  67 = 1;  68 = 0;  69 = 1;  70 = 1;  71 = 2;  72 = 1;  73 = 1;
    74 = 2;  75 = 2;  76 = 1;  77 = 1;  78 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
